$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.002.81'
$ws.Range("E2").Value = '  -0.22%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.827.37'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.009'
$ws.Range("E4").Value = '  -0.28%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '311.66'
$ws.Range("E5").Value = '  -0.34%  '
$ws.Range("E6").Value = '  -0.22%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4656'
$ws.Range("E7").Value = '  -1.05%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3708'
$ws.Range("E8").Value = '  +1.56%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07373'
$ws.Range("E9").Value = '  -0.34%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8738'
$ws.Range("E10").Value = '  -0.69%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '19.99'
$ws.Range("E11").Value = '  -1.63%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07829'
$ws.Range("E12").Value = '  +6.97%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.825.51'
$ws.Range("E13").Value = '  -5.03%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.355'
$ws.Range("E14").Value = '  -0.36%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.563'
$ws.Range("E15").Value = '  +0.72%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '91.95'
$ws.Range("E16").Value = '  -1.38%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.010'
$ws.Range("E17").Value = '  +0.08%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008894'
$ws.Range("E18").Value = '  +2.23%  '
$ws.Range("E19").Value = '  -0.35%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.69'
$ws.Range("E20").Value = '  +0.45%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '26.886.74'
$ws.Range("E21").Value = '  -2.74%  '
$ws.Range("E22").Value = '  -1.62%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.57'
$ws.Range("E23").Value = '  -0.22%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.068.72'
$ws.Range("E24").Value = '  -1.42%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '152.58'
$ws.Range("E25").Value = '  +0.54%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.831'
$ws.Range("E26").Value = '  -2.95%  '
$ws.Range("E27").Value = '  -1.42%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.099'
$ws.Range("E28").Value = '  -1.97%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.124'
$ws.Range("E29").Value = '  -0.87%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '115.50'
$ws.Range("E30").Value = '  -0.47%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08865'
$ws.Range("E31").Value = '  -0.77%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.972'
$ws.Range("E32").Value = '  +0.75%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.7291'
$ws.Range("E33").Value = '  -1.49%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.443'
$ws.Range("E34").Value = '  -1.41%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.136'
$ws.Range("E35").Value = '  -2.43%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.517'
$ws.Range("E36").Value = '  +4.82%  '
$ws.Range("E37").Value = '  +0.76%  '
$ws.Range("E38").Value = '  -1.12%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05244'
$ws.Range("E39").Value = '  -0.86%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '7.265'
$ws.Range("E40").Value = '  +0.93%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.928'
$ws.Range("E41").Value = '  -0.05%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.5192'
$ws.Range("E42").Value = '  -0.94%  '
$ws.Range("E43").Value = '  -14.35%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.1627'
$ws.Range("E44").Value = '  -0.79%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '8.226'
$ws.Range("E45").Value = '  -1.89%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4852'
$ws.Range("E46").Value = '  -0.28%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '10.26'
$ws.Range("E47").Value = '  -1.32%  '
$ws.Range("E48").Value = '  -0.21%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '102.80'
$ws.Range("E49").Value = '  -1.31%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.624'
$ws.Range("E50").Value = '  -1.72%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06216'
$ws.Range("E51").Value = '  -1.20%  '
